$wb = $excel.ActiveWorkbook

# --- "0D" sheet: B6 was "=1/1000" -> becomes a plain literal value 1 ---
$ws0D = $wb.Worksheets.Item("0D")
$ws0D.Range("B6").Value = 1
# Update the sheet's selection (no longer the active tab after this edit)
$ws0D.Range("B7").Select()

# --- "Vecteurs" sheet: several Emissions_scope_2_3 values divided by 1000 ---
$wsVec = $wb.Worksheets.Item("Vecteurs")
$wsVec.Range("C2").Value = 0.06
$wsVec.Range("C3").Value = 0.006
$wsVec.Range("C6").Value = 0.2
$wsVec.Range("C8").Value = 0.03
$wsVec.Range("C9").Value = 0.006
$wsVec.Range("C12").Value = 0.06

# "Vecteurs" becomes the active sheet/tab, with a new selection
$wsVec.Activate()
$wsVec.Range("G12").Select()
